$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.890.92'
$ws.Range("E2").Value = '  +1.66%  '
$ws.Range("D3").Value = '3.129.01'
$ws.Range("E3").Value = '  +2.06%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = '''577.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").Value = '''179.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.05%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.119.25'
$ws.Range("E8").Value = '  +1.81%  '
$ws.Range("E9").Value = '  +1.30%  '
$ws.Range("E10").Value = '  +2.57%  '
$ws.Range("E11").Value = '  +1.90%  '
$ws.Range("D12").Value = '''0.470'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("E13").Value = '  +1.12%  '
$ws.Range("D14").Value = '''36.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.28%  '
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '3.644.59'
$ws.Range("E16").Value = '  +1.98%  '
$ws.Range("D17").Value = '67.689.58'
$ws.Range("E17").Value = '  +1.46%  '
$ws.Range("D18").Value = '''7.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("D19").Value = '3.121.94'
$ws.Range("E19").Value = '  +1.96%  '
$ws.Range("D20").Value = '''16.43'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.98%  '
$ws.Range("D21").Value = '''487.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("B22").Value = 'Polygon'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D22").Value = '''0.693'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '''7.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("D24").Value = '''83.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.46%  '
$ws.Range("D25").Value = '''2.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.38%  '
$ws.Range("D26").Value = '''12.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.17%  '
$ws.Range("E27").Value = '  +4.88%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '''8.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.93%  '
$ws.Range("D30").Value = '''2.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.89%  '
$ws.Range("E31").Value = '  +0.73%  '
$ws.Range("D32").Value = '''28.24'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.70%  '
$ws.Range("E33").Value = '  +1.20%  '
$ws.Range("E34").Value = '  +4.33%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = '''48.85'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.85%  '
$ws.Range("D37").Value = '''0.955'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.78%  '
$ws.Range("D38").Value = '''5.62'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.36%  '
$ws.Range("D39").Value = '''0.324'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.16%  '
$ws.Range("D40").Value = '''2.04'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.23%  '
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("E42").Value = '  +1.36%  '
$ws.Range("D43").Value = '''8.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.53%  '
$ws.Range("E44").Value = '  +8.71%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").Value = '''385.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.35%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.779.99'
$ws.Range("E46").Value = '  +0.94%  '
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("D48").Value = '''26.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.98%  '
$ws.Range("D49").Value = '''136.15'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = '''2.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.41%  '
